$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('!!_Table of contents')
$ws.Unprotect()
$ws.Range('A1').Value2 = '!!!ObjTables objTablesVersion=''1.0.0'' date=''2020-05-29 00:18:55'''
$ws.Range('A2').Value2 = '!!ObjTables type=''TableOfContents'' tableFormat=''row'' description=''Table of contents'' date=''2020-05-29 00:18:55'' objTablesVersion=''1.0.0'''
$ws.Protect()

$ws = $wb.Worksheets.Item('!!_Schema')
$ws.Unprotect()
$ws.Range('A1').Value2 = '!!ObjTables type=''Schema'' tableFormat=''row'' description=''Table/model and column/attribute definitions'' date=''2020-05-29 00:18:55'' objTablesVersion=''1.0.0'''
$ws.Range('D4').Value2 = 'Slug(r''^(?!(^|\b)(\d+(\.\d*)?(\b|$))|(\.\d+$)|(0[x][0-9a-f]+(\b|$))|([0-9]+e[0-9]+(\b|$)))[a-z0-9_]+$'', flags=2, primary=True, unique=True)'
$ws.Range('E4').Value2 = 'Id'
$ws.Range('E5').Value2 = 'Identifiers'
$ws.Range('E6').Value2 = 'IsConstant'
$ws.Range('E7').Value2 = 'Model'
$ws.Range('E8').Value2 = 'Name'
$ws.Range('D10').Value2 = 'Slug(r''^(?!(^|\b)(\d+(\.\d*)?(\b|$))|(\.\d+$)|(0[x][0-9a-f]+(\b|$))|([0-9]+e[0-9]+(\b|$)))[a-z0-9_]+$'', flags=2, primary=True, unique=True)'
$ws.Range('E10').Value2 = 'Id'
$ws.Range('E11').Value2 = 'Name'
$ws.Range('E13').Value2 = 'Equation'
$ws.Range('E14').Value2 = 'Gene'
$ws.Range('D15').Value2 = 'Slug(r''^(?!(^|\b)(\d+(\.\d*)?(\b|$))|(\.\d+$)|(0[x][0-9a-f]+(\b|$))|([0-9]+e[0-9]+(\b|$)))[a-z0-9_]+$'', flags=2, primary=True, unique=True)'
$ws.Range('E15').Value2 = 'Id'
$ws.Range('E16').Value2 = 'Identifiers'
$ws.Range('E17').Value2 = 'IsReversible'
$ws.Range('E18').Value2 = 'Model'
$ws.Range('E19').Value2 = 'Name'
$ws.Protect()

$ws = $wb.Worksheets.Item('!!Compound')
$ws.Unprotect()
$ws.Range('A1').Value2 = '!!ObjTables type=''Data'' tableFormat=''row'' class=''Compound'' name=''Compound'' description=''Compound'' date=''2020-03-10 22:56:34'' objTablesVersion=''1.0.0'''
$ws.Protect()

$ws = $wb.Worksheets.Item('!!Model')
$ws.Unprotect()
$ws.Range('A1').Value2 = '!!ObjTables type=''Data'' tableFormat=''column'' class=''Model'' name=''Model'' description=''Model'' date=''2020-03-10 22:56:35'' objTablesVersion=''1.0.0'''
$ws.Protect()

$ws = $wb.Worksheets.Item('!!Reaction')
$ws.Unprotect()
$ws.Range('A1').Value2 = '!!ObjTables type=''Data'' tableFormat=''row'' class=''Reaction'' name=''Reaction'' description=''Reaction'' date=''2020-03-10 22:56:35'' objTablesVersion=''1.0.0'''
$ws.Protect()
